$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 - "Odd Even Linked List" (leetcode 328)
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 328
$ws.Range("C17").Value = "给定一个单链表，把所有的奇数节点和偶数节点分别排在一起。请注意，这里的奇数节点和偶数节点指的是节点编号的奇偶性，而不是节点的值的奇偶性"
$ws.Range("D17").Value = "1 初始化奇偶位置头指针【迭代期间不动】以及奇偶位置移动指针【负责迭代】
2 奇指针指向偶指针的next节点
3 奇指针向后迭代【偶指针之后的节点】
4 偶指针指向奇指针之后的节点
5 偶指针向后迭代【奇指针之后的节点】
6 判断偶指针以及偶指针的next是否是null【注意】
7 奇数个节点的链表，最后的状态是even_cur链表尾部的null；odd_cur指向最后一个节点
7 偶数个节点的链表，最后的状态是even_cur链表最后一个节点；odd_cur指向链表倒数第二个节点
8 奇链表尾部的next【odd_cur】指向偶链表的头部
9 返回奇链表的头部"
$ws.Range("E17").Value = "链表添加
保存节点"
$ws.Range("F17").Value = "O(N), N是元素个数"
$ws.Range("G17").Value = "O(1)"
$ws.Rows.Item(17).RowHeight = 260

# Row 18 - "Add Two Numbers" (leetcode 445)
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 445
$ws.Range("C18").Value = "给你两个 非空 链表来代表两个非负整数。数字最高位位于链表开始位置。它们的每个节点只存储一位数字。将这两数相加会返回一个新的链表。"
$ws.Range("D18").Value = "1 两个链表分别入栈，入栈完毕
2 取出栈顶元素，相加，获取相加之和【一位数字】以及进位
3 循环：根据创建节点tmp，solder指向tmp，下一个节点要添加到solder之后
4 考虑最高位有进位的情况。例如9，9"
$ws.Range("E18").Value = "栈
链表反转
链表添加"
$ws.Range("F18").Value = "O(M+N),M,N是两个链表的元素个数"
$ws.Range("G18").Value = "O(1)"
$ws.Rows.Item(18).RowHeight = 120

$ws.Range("C21").Select() | Out-Null
